# Update Name of Algo
# Applies updated RandomForest imputation results for columns D and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Cell = "E3"; Value = 16.1946 },
    @{ Cell = "D12"; Value = -6.7544 },
    @{ Cell = "E14"; Value = 16.51090000000001 },
    @{ Cell = "E26"; Value = 16.54749999999999 },
    @{ Cell = "D27"; Value = -8.420700000000004 },
    @{ Cell = "E31"; Value = 16.5928 },
    @{ Cell = "D32"; Value = -8.584199999999994 },
    @{ Cell = "E35"; Value = 16.80089999999999 },
    @{ Cell = "D36"; Value = -8.026499999999999 },
    @{ Cell = "E37"; Value = 16.7296 },
    @{ Cell = "D38"; Value = -7.780799999999996 },
    @{ Cell = "E45"; Value = 16.5239 },
    @{ Cell = "D46"; Value = -8.381599999999999 },
    @{ Cell = "E52"; Value = 17.23350000000001 },
    @{ Cell = "D54"; Value = -8.348300000000002 },
    @{ Cell = "D55"; Value = -8.294299999999994 },
    @{ Cell = "D56"; Value = -8.072499999999996 },
    @{ Cell = "E57"; Value = 16.6046 },
    @{ Cell = "D67"; Value = -6.344099999999997 },
    @{ Cell = "D69"; Value = -7.072699999999995 },
    @{ Cell = "D72"; Value = -7.346799999999999 },
    @{ Cell = "E81"; Value = 16.0537 },
    @{ Cell = "D83"; Value = -8.962600000000002 },
    @{ Cell = "E83"; Value = 16.5193 },
    @{ Cell = "D86"; Value = -7.618599999999994 },
    @{ Cell = "D91"; Value = -6.691199999999998 },
    @{ Cell = "D93"; Value = -6.6146 },
    @{ Cell = "D99"; Value = -7.514100000000003 },
    @{ Cell = "E100"; Value = 16.4224 },
    @{ Cell = "E102"; Value = 16.7993 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$wb.Save()
